$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-08"

# Update the "April (through 04-07)" label to "April (through 04-08)"
$ws.Range("A5").Value = "April (through 04-08)"

# Update April row (row 5) values
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 19
$ws.Range("I5").Value = 27

# Update Total row (row 6) values
$ws.Range("B6").Value = 71
$ws.Range("C6").Value = 134
$ws.Range("D6").Value = 203
$ws.Range("E6").Value = 211
$ws.Range("F6").Value = 124
$ws.Range("G6").Value = 213
$ws.Range("H6").Value = 442
$ws.Range("I6").Value = 461
